$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 311; this shifts the existing rows 311-317
# down to 312-318 (and extends the sheet dimension accordingly).
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new weekly price record.
$ws.Cells.Item(311, 1).Value  = 5
$ws.Cells.Item(311, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(311, 3).Value  = "Maule"
$ws.Cells.Item(311, 4).Value  = 44448
$ws.Cells.Item(311, 5).Value  = 7
$ws.Cells.Item(311, 6).Value  = 100112020
$ws.Cells.Item(311, 7).Value  = "Tomate"
$ws.Cells.Item(311, 8).Value  = "Larga vida"
$ws.Cells.Item(311, 9).Value  = "Primera"
$ws.Cells.Item(311, 10).Value = 2500
$ws.Cells.Item(311, 11).Value = 10000
$ws.Cells.Item(311, 12).Value = 10000
$ws.Cells.Item(311, 13).Value = 10000
$ws.Cells.Item(311, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(311, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(311, 16).Value = 1000
$ws.Cells.Item(311, 17).Value = 10
$ws.Cells.Item(311, 18).Value = "Hortaliza"
